$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-33 from 2023-09-23
# (serial 45192) to 2023-10-03 (serial 45202), a 10 day shift, while keeping
# the existing date formatting/style of each cell.
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}
